$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.333.34'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.269.19'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '318.03'
$ws.Range("E5").Value = '  -1.36%  '
$ws.Range("D6").Value = '100.24'
$ws.Range("E6").Value = '  -5.02%  '
$ws.Range("E7").Value = '  -2.21%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -5.73%  '
$ws.Range("D10").Value = '36.46'
$ws.Range("E10").Value = '  -6.00%  '
$ws.Range("D11").Value = '0.0830'
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").Value = '7.45'
$ws.Range("E12").Value = '  -5.66%  '
$ws.Range("E13").Value = '  -2.49%  '
$ws.Range("D14").Value = '2.614.76'
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = '0.851'
$ws.Range("E15").Value = '  -4.03%  '
$ws.Range("D16").Value = '2.263.87'
$ws.Range("E17").Value = '  -3.48%  '
$ws.Range("D18").Value = '44.271.06'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '13.27'
$ws.Range("E19").Value = '  -4.94%  '
$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("E20").Value = '  -2.09%  '
$ws.Range("D21").Value = '6.41'
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").Value = '65.90'
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").Value = '241.06'
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("E24").Value = '  -6.17%  '
$ws.Range("D25").Value = '2.06'
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").Value = '10.26'
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").Value = '38.69'
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  -3.45%  '
$ws.Range("D30").Value = '6.08'
$ws.Range("E30").Value = '  -6.78%  '
$ws.Range("D31").Value = '20.24'
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("D32").Value = '155.61'
$ws.Range("E32").Value = '  -4.15%  '
$ws.Range("D33").Value = '0.0846'
$ws.Range("E33").Value = '  -4.69%  '
$ws.Range("D34").Value = '3.48'
$ws.Range("E34").Value = '  +10.33%  '
$ws.Range("E35").Value = '  -4.00%  '
$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  -4.62%  '
$ws.Range("E37").Value = '  -4.74%  '
$ws.Range("E38").Value = '  -2.65%  '
$ws.Range("D39").Value = '15.55'
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").Value = '3.57'
$ws.Range("E40").Value = '  -10.05%  '
$ws.Range("D41").Value = '3.96'
$ws.Range("D42").Value = '0.0311'
$ws.Range("E42").Value = '  -5.86%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = '1.740.80'
$ws.Range("E44").Value = '  -2.51%  '
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("E46").Value = '  -6.16%  '
$ws.Range("D47").Value = '5.24'
$ws.Range("D48").Value = '102.78'
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D49").Value = '57.17'
$ws.Range("E49").Value = '  -6.24%  '
$ws.Range("E50").Value = '  -5.16%  '
$ws.Range("D51").Value = '71.53'
$ws.Range("E51").Value = '  -5.44%  '
